$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 978.2778
$ws.Range("I15").Value = 978.2778
$ws.Range("K15").Value = 2934.8334
$ws.Range("M15").Value = -2765.8334
# Row 40
$ws.Range("H40").Value = 4173
$ws.Range("I40").Value = 2270.2856
$ws.Range("J40").Value = 4874
$ws.Range("K40").Value = 2270.2856
$ws.Range("L40").Value = 4874
$ws.Range("M40").Value = -2095.2856
$ws.Range("N40").Value = -5224
# Row 57
$ws.Range("H57").Value = 69860
$ws.Range("J57").Value = 69860
$ws.Range("L57").Value = 209580
$ws.Range("N57").Value = -210578
# Row 92
$ws.Range("H92").Value = 1783.0741
$ws.Range("I92").Value = 1035.9524
$ws.Range("K92").Value = 1035.9524
$ws.Range("M92").Value = 212.0476000000001
# Row 98
$ws.Range("H98").Value = 1317.3334
$ws.Range("I98").Value = 1259.0555
$ws.Range("J98").Value = 1667
$ws.Range("K98").Value = 1259.0555
$ws.Range("L98").Value = 1667
$ws.Range("M98").Value = 238.9445000000001
$ws.Range("N98").Value = -4663
# Row 106
$ws.Range("H106").Value = 31251852
$ws.Range("I106").Value = 33335174
$ws.Range("K106").Value = 33335174
$ws.Range("M106").Value = -33334543
# Row 122
$ws.Range("H122").Value = 1317.3334
$ws.Range("I122").Value = 1259.0555
$ws.Range("J122").Value = 1667
$ws.Range("K122").Value = 3777.1665
$ws.Range("L122").Value = 5001
$ws.Range("M122").Value = -1327.1665
$ws.Range("N122").Value = -9901
# Row 132
$ws.Range("H132").Value = 5812.2383
$ws.Range("I132").Value = 6721.2354
$ws.Range("J132").Value = 1949
$ws.Range("K132").Value = 20163.7062
$ws.Range("L132").Value = 5847
$ws.Range("M132").Value = -17633.7062
$ws.Range("N132").Value = -10907
# Row 137
$ws.Range("H137").Value = 70926.08
$ws.Range("I137").Value = 138723.69
$ws.Range("K137").Value = 416171.07
$ws.Range("M137").Value = -413621.07
# Row 138
$ws.Range("H138").Value = 3363.2952
$ws.Range("I138").Value = 2614.4285
$ws.Range("J138").Value = 3586.3618
$ws.Range("K138").Value = 7843.2855
$ws.Range("L138").Value = 10759.0854
$ws.Range("M138").Value = -2703.2855
$ws.Range("N138").Value = -21039.0854

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6842.644
$ws.Range("I32").Value = 5854.25
$ws.Range("J32").Value = 20284.8
$ws.Range("K32").Value = 5854.25
$ws.Range("L32").Value = 20284.8
$ws.Range("M32").Value = -5567.25
$ws.Range("N32").Value = -20858.8
# Row 61
$ws.Range("H61").Value = 5034.7856
$ws.Range("I61").Value = 5046.091
$ws.Range("J61").Value = 4993.3335
$ws.Range("K61").Value = 5046.091
$ws.Range("L61").Value = 4993.3335
$ws.Range("M61").Value = -4834.091
$ws.Range("N61").Value = -5417.3335
# Row 112
$ws.Range("H112").Value = 32846.75
$ws.Range("J112").Value = 32846.75
$ws.Range("L112").Value = 32846.75
$ws.Range("N112").Value = -35800.75
# Row 132
$ws.Range("H132").Value = 3336.3684
$ws.Range("I132").Value = 2111.4443
$ws.Range("J132").Value = 4438.8
$ws.Range("K132").Value = 6334.3329
$ws.Range("L132").Value = 13316.4
$ws.Range("M132").Value = -3804.3329
$ws.Range("N132").Value = -18376.4
# Row 136
$ws.Range("H136").Value = 5034.7856
$ws.Range("I136").Value = 5046.091
$ws.Range("J136").Value = 4993.3335
$ws.Range("K136").Value = 15138.273
$ws.Range("L136").Value = 14980.0005
$ws.Range("M136").Value = -12588.273
$ws.Range("N136").Value = -20080.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3473870.2
$ws.Range("I105").Value = 3473870.2
$ws.Range("K105").Value = 3473870.2
$ws.Range("M105").Value = -3472123.2
# Row 107
$ws.Range("H107").Value = 2382640.2
$ws.Range("I107").Value = 3402725.8
$ws.Range("J107").Value = 2440.4443
$ws.Range("K107").Value = 3402725.8
$ws.Range("L107").Value = 2440.4443
$ws.Range("M107").Value = -3400805.8
$ws.Range("N107").Value = -6280.4443
# Row 112
$ws.Range("H112").Value = 88746.25
$ws.Range("J112").Value = 88746.25
$ws.Range("L112").Value = 88746.25
$ws.Range("N112").Value = -91700.25
# Row 134
$ws.Range("H134").Value = 7917.1875
$ws.Range("I134").Value = 1548.2142
$ws.Range("K134").Value = 4644.642599999999
$ws.Range("M134").Value = -2109.642599999999

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 96.052635
$ws.Range("I7").Value = 35.5
$ws.Range("K7").Value = 35.5
$ws.Range("M7").Value = 77.5
# Row 28
$ws.Range("H28").Value = 19760.75
$ws.Range("J28").Value = 19760.75
$ws.Range("L28").Value = 19760.75
$ws.Range("N28").Value = -20250.75
# Row 107
$ws.Range("H107").Value = 1201.3019
$ws.Range("J107").Value = 847.6667
$ws.Range("L107").Value = 847.6667
$ws.Range("N107").Value = -4687.6667
# Row 134
$ws.Range("H134").Value = 3825.0588
$ws.Range("I134").Value = 2419.6667
$ws.Range("K134").Value = 7259.000100000001
$ws.Range("M134").Value = -4724.000100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 100057620
$ws.Range("I32").Value = 67999.664
$ws.Range("J32").Value = 160051380
$ws.Range("K32").Value = 203998.992
$ws.Range("L32").Value = 480154140
$ws.Range("M32").Value = -203715.992
$ws.Range("N32").Value = -480154706
# Row 122
$ws.Range("H122").Value = 877.4
$ws.Range("I122").Value = 999.8889
$ws.Range("J122").Value = 808.5
$ws.Range("K122").Value = 8999.000100000001
$ws.Range("L122").Value = 7276.5
$ws.Range("M122").Value = -6549.000100000001
$ws.Range("N122").Value = -12176.5

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 59453.117
$ws.Range("J107").Value = 646
$ws.Range("L107").Value = 646
$ws.Range("N107").Value = -4486
# Row 111
$ws.Range("H111").Value = 28539
$ws.Range("J111").Value = 28539
$ws.Range("L111").Value = 28539
$ws.Range("N111").Value = -34673
# Row 126
$ws.Range("H126").Value = 4957590.5
$ws.Range("I126").Value = 4135807
$ws.Range("K126").Value = 12407421
$ws.Range("M126").Value = -12404951
# Row 132
$ws.Range("H132").Value = 3355.1428
$ws.Range("I132").Value = 3221.75
$ws.Range("J132").Value = 3688.625
$ws.Range("K132").Value = 9665.25
$ws.Range("L132").Value = 11065.875
$ws.Range("M132").Value = -7135.25
$ws.Range("N132").Value = -16125.875
# Row 136
$ws.Range("H136").Value = 12165.156
$ws.Range("J136").Value = 12165.156
$ws.Range("L136").Value = 36495.468
$ws.Range("N136").Value = -41595.468

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5326.7427
$ws.Range("I7").Value = 3703.6086
$ws.Range("J7").Value = 8437.75
$ws.Range("K7").Value = 3703.6086
$ws.Range("L7").Value = 8437.75
$ws.Range("M7").Value = -3591.6086
$ws.Range("N7").Value = -8661.75
# Row 61
$ws.Range("H61").Value = 5051463
$ws.Range("I61").Value = 5848845
$ws.Range("K61").Value = 5848845
$ws.Range("M61").Value = -5848643
# Row 110
$ws.Range("H110").Value = 23000
$ws.Range("J110").Value = 23000
$ws.Range("L110").Value = 23000
$ws.Range("N110").Value = -31180
# Row 113
$ws.Range("H113").Value = 5051463
$ws.Range("I113").Value = 5848845
$ws.Range("K113").Value = 5848845
$ws.Range("M113").Value = -5846675
# Row 122
$ws.Range("H122").Value = 6136
$ws.Range("I122").Value = 3942.1428
$ws.Range("J122").Value = 7842.3335
$ws.Range("K122").Value = 11826.4284
$ws.Range("L122").Value = 23527.0005
$ws.Range("M122").Value = -9376.428400000001
$ws.Range("N122").Value = -28427.0005
# Row 126
$ws.Range("H126").Value = 5326.7427
$ws.Range("I126").Value = 3703.6086
$ws.Range("J126").Value = 8437.75
$ws.Range("K126").Value = 11110.8258
$ws.Range("L126").Value = 25313.25
$ws.Range("M126").Value = -8640.825800000001
$ws.Range("N126").Value = -30253.25
# Row 132
$ws.Range("H132").Value = 7258.5103
$ws.Range("I132").Value = 7330.452
$ws.Range("K132").Value = 21991.356
$ws.Range("M132").Value = -19461.356

$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 11416.167
$ws.Range("I26").Value = 7500
$ws.Range("J26").Value = 12199.4
$ws.Range("K26").Value = 7500
$ws.Range("L26").Value = 12199.4
$ws.Range("M26").Value = -7207
$ws.Range("N26").Value = -12785.4
# Row 107
$ws.Range("H107").Value = 50001104
$ws.Range("I107").Value = 111112020
$ws.Range("K107").Value = 333336060
$ws.Range("M107").Value = -333334140
# Row 110
$ws.Range("H110").Value = 38314.668
$ws.Range("J110").Value = 38314.668
$ws.Range("L110").Value = 38314.668
$ws.Range("N110").Value = -46494.668
# Row 113
$ws.Range("H113").Value = 1019.8077
$ws.Range("I113").Value = 187.17647
$ws.Range("K113").Value = 561.52941
$ws.Range("M113").Value = 1608.47059
# Row 121
$ws.Range("H121").Value = 79207.5
$ws.Range("J121").Value = 79207.5
$ws.Range("L121").Value = 79207.5
$ws.Range("N121").Value = -82701.5
# Row 132
$ws.Range("H132").Value = 17738524
$ws.Range("I132").Value = 21743008
$ws.Range("J132").Value = 992502.8
$ws.Range("K132").Value = 65229024
$ws.Range("L132").Value = 2977508.4
$ws.Range("M132").Value = -65226494
$ws.Range("N132").Value = -2982568.4
# Row 137
$ws.Range("H137").Value = 84333.164
$ws.Range("J137").Value = 84333.164
$ws.Range("L137").Value = 84333.164
$ws.Range("N137").Value = -94533.164
